# Word COM-interop script implementing the diff:
#   - "Configure .bashrc file to add path" -> "Configure .bashrc for PATH variables"
#   - Insert four new sub-bullets (numId 1010 / ilvl 1, style "Compact") right after
#     that bullet and before "Demonstrate installation success with --version commands",
#     each reading "Add: <export ...>" with the export command in the VerbatimChar
#     character style.

$d = $word.ActiveDocument

# 1) Fix the wording on the "Configure .bashrc ..." bullet.
$r1 = $d.Content
$null = $r1.Find.Execute("file to add path", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "for PATH variables", 2)

# 2) Insert the four new "Add: export ..." bullets in front of the
#    "Demonstrate installation success ..." bullet, which already carries the
#    numId 1010 / ilvl 1 "Compact" list formatting that the new bullets need.
#    Using Find/Replace with "^p" paragraph marks means every new paragraph
#    inherits that exact formatting automatically.
$lines = @(
    "Add: export JAVA_HOME=/usr",
    "Add: export HADOOP_HOME=/home/vagrant/hadoop-2.9.2",
    "Add: export SPARK_HOME=/home/vagrant/spark",
    "Add: export HADOOP_CLASSPATH=/usr/lib/jvm/java-8-openjdk-amd64/lib/tools.jar"
)

$replacement = ($lines -join "^p") + "^pDemonstrate installation success"

$r2 = $d.Content
$null = $r2.Find.Execute("Demonstrate installation success", $true, $false, $false, $false, $false, `
                          $true, 1, $false, $replacement, 2)

# 3) Re-find each new bullet's "export ..." text and apply the VerbatimChar
#    character style to it, leaving the leading "Add: " text plain.
foreach ($line in $lines) {
    $verbatim = $line.Substring(5)   # strip the leading "Add: "

    $r3 = $d.Content
    $found = $r3.Find.Execute($verbatim, $false, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
    if ($found) {
        $r3.Style = "VerbatimChar"
    }
}

Write-Host "Edit applied"
